$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.09"
$ws.Range("E2").Value = "'-4.59%"
$ws.Range("E3").Value = "'-0.95%"
$ws.Range("D4").Value = "'5.124"
$ws.Range("E4").Value = "'-4.16%"
$ws.Range("D5").Value = "'0.07495"
$ws.Range("E5").Value = "'-1.14%"
$ws.Range("D6").Value = "'7.739"
$ws.Range("E6").Value = "'-1.47%"
$ws.Range("D7").Value = "'1.708"
$ws.Range("E7").Value = "'5.49%"
$ws.Range("D8").Value = "'3.800"
$ws.Range("E8").Value = "'2.28%"
$ws.Range("D9").Value = "'0.9322"
$ws.Range("E9").Value = "'1.68%"
$ws.Range("D10").Value = "'0.1695"
$ws.Range("E10").Value = "'-1.74%"
$ws.Range("D11").Value = "'0.07185"
$ws.Range("E11").Value = "'-6.90%"
$ws.Range("D12").Value = "'0.07920"
$ws.Range("E12").Value = "'-3.50%"
$ws.Range("D13").Value = "'0.03022"
$ws.Range("E13").Value = "'0.29%"
$ws.Range("D14").Value = "'0.09896"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("D15").Value = "'0.001512"
$ws.Range("E15").Value = "'-0.76%"
$ws.Range("D16").Value = "'0.006303"
$ws.Range("E16").Value = "'-2.44%"
$ws.Range("D17").Value = "'3.451"
$ws.Range("E17").Value = "'-1.15%"
$ws.Range("D18").Value = "'2.228"
$ws.Range("E18").Value = "'-0.61%"
$ws.Range("E19").Value = "'-0.86%"
$ws.Range("D20").Value = "'0.1328"
$ws.Range("D21").Value = "'4.555"
$ws.Range("E21").Value = "'8.56%"
$ws.Range("D22").Value = "'0.04645"
$ws.Range("E22").Value = "'1.93%"
$ws.Range("D23").Value = "'0.1560"
$ws.Range("E23").Value = "'-3.69%"
$ws.Range("E24").Value = "'-0.18%"
$ws.Range("D25").Value = "'0.004427"
$ws.Range("E25").Value = "'-1.48%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("D27").Value = "'0.0001876"
$ws.Range("E27").Value = "'7.92%"
$ws.Range("D39").Value = "'0.01672"
$ws.Range("E39").Value = "'-1.15%"
$ws.Range("D40").Value = "'0.04461"
$ws.Range("E40").Value = "'-3.14%"
$ws.Range("D41").Value = "'0.007078"
$ws.Range("E41").Value = "'-1.86%"
$ws.Range("D42").Value = "'0.1327"
$ws.Range("E42").Value = "'-2.86%"
$ws.Range("D43").Value = "'0.002061"
$ws.Range("E43").Value = "'-8.84%"
$ws.Range("D44").Value = "'0.01136"
$ws.Range("E44").Value = "'-19.01%"
$ws.Range("D45").Value = "'0.00006002"
$ws.Range("E45").Value = "'-2.23%"
$ws.Range("D46").Value = "'1.930"
$ws.Range("E46").Value = "'1.96%"
$ws.Range("E47").Value = "'-0.13%"
